$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strings for the added rows (227-230)
$fStr227 = "0:17`n7. Armazenamento de imagens usando Amazon S3`n88. Salvando URL da imagem em Cliente`nabordagem provisória - salvar a URL da imagem em Cliente"
$fStr228 = "`n1:00`n7. Armazenamento de imagens usando Amazon S3`n88. Salvando URL da imagem em Cliente`nretirando inserção de imagem (permissão do endpoint ""/clientes/picture"")que foi inserida provisóriamente para testes e exigir login/autenticação para que seja inserido imagens no banco"
$fStr229 = "3:58`n7. Armazenamento de imagens usando Amazon S3`n88. Salvando URL da imagem em Cliente`nno video o professor usa o ""repo.FindOne(user.getId())"" ao instanciar o cliente mas este metodo gera erro ... para consertar usa-se: ""Cliente cli = find(user.getId());"""
$fStr230 = "6:57`n7. Armazenamento de imagens usando Amazon S3`n88. Salvando URL da imagem em Cliente`npara enviar a imagem de acordo com o usuario autenticado, é necessário logar, depois pegar o token no header, acessar o endpoint, incluir a imagem junto do token para que funcione"
$eStr = "Salvando URL da imagem em Cliente"
$gBlank = "`n`n`n`n`n"

# Extend the table (Tabela1) by 4 rows - this grows the table/autofilter ref
# from B1:G226 to B1:G230
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Duplicate the formatting of the last previously-existing row (226) down into
# the 4 new rows
$ws.Range("B226:G226").Copy()
$ws.Range("B227:G230").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Set B/C/D (section/name/lesson number) for all 4 new rows first
$ws.Cells.Item(227, 2).Value = 7
$ws.Cells.Item(227, 3).Value = "Armazenamento de imagens usando Amazon S3"
$ws.Cells.Item(227, 4).Value = 88

$ws.Cells.Item(228, 2).Value = 7
$ws.Cells.Item(228, 3).Value = "Armazenamento de imagens usando Amazon S3"
$ws.Cells.Item(228, 4).Value = 88

$ws.Cells.Item(229, 2).Value = 7
$ws.Cells.Item(229, 3).Value = "Armazenamento de imagens usando Amazon S3"
$ws.Cells.Item(229, 4).Value = 88

$ws.Cells.Item(230, 2).Value = 7
$ws.Cells.Item(230, 3).Value = "Armazenamento de imagens usando Amazon S3"
$ws.Cells.Item(230, 4).Value = 88

# Set column F (abordagem da aula) for all 4 rows first, so the new shared
# strings are created in this exact order (matches uniqueCount indices 331-334)
$ws.Cells.Item(227, 6).Value = $fStr227
$ws.Cells.Item(228, 6).Value = $fStr228
$ws.Cells.Item(229, 6).Value = $fStr229
$ws.Cells.Item(230, 6).Value = $fStr230

# Set column E (nome aula) afterwards - creates the last new shared string (335)
# on first use, then the other three rows simply reuse it
$ws.Cells.Item(227, 5).Value = $eStr
$ws.Cells.Item(228, 5).Value = $eStr
$ws.Cells.Item(229, 5).Value = $eStr
$ws.Cells.Item(230, 5).Value = $eStr

# Column G (aprendido) only has content on the first new row, reusing an
# already-existing shared string made up of blank lines
$ws.Cells.Item(227, 7).Value = $gBlank

# Row heights matching target layout
$ws.Rows.Item(227).RowHeight = 90
$ws.Rows.Item(228).RowHeight = 105
$ws.Rows.Item(229).RowHeight = 90
$ws.Rows.Item(230).RowHeight = 90

# Update view: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 225
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E231").Select()
